$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("2:4").Insert()
$ws.Range("A2").Value = "BaseStateChromium"
$ws.Range("B2").Value = "Started"
$ws.Range("C2").Value = "Not Started"
$ws.Range("A3").Value = "BaseStatePixel9Pro_API35"
$ws.Range("B3").Value = "Not Started"
$ws.Range("C3").Value = "Started"
$ws.Range("A4").Value = "AUT"
$ws.Range("B4").Value = "Chromium"
$ws.Range("C4").Value = "Pixel9Pro_API35"
$ws.Range("B14").Select() | Out-Null
